$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.404.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.938.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7702"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "249.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.04"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3205"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.30%  "
$ws.Range("E10").Value = "  -2.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7856"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08009"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.939.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.389"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "95.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("E16").Value = "  -3.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.405.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "257.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008030"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.852"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.190.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.775"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.627"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1342"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.301"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.367"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.446"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.165"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05199"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.284"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7525"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.781"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.30%  "
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.805"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.468"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4529"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.983"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8365"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.812"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.539"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "987.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.40%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.38%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4169"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.30%  "
